# Cycle monitoring workbook update
# - Row 20: record the low Ammonium phosphate / Nitrite readings (G20, H20)
# - Row 21: fill in the rest of day-20 readings + a note about dosing ammonium chloride
# - Rows 22-26: new daily readings (days 21-25)
# - Rows 27-31: new placeholder days (26-30), same bare pattern as the old
#   trailing rows (date / cycle day / instructions only)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$instructions = "Measure pH, Ammonia, Nitrite and Nitrate and record values"

# ---- Row 20: add the two measurements that were already pending ----
$ws.Range("G19:H19").Copy()
$ws.Range("G20:H20").PasteSpecial(-4122)
$ws.Range("G20").Value = 0.05
$ws.Range("H20").Value = 0.1

# ---- Row 21: bring the rest of the row up to the same shape as row 19,
#      plus a note in column L ----
$ws.Range("D19:I19").Copy()
$ws.Range("D21:I21").PasteSpecial(-4122)
$ws.Range("L4").Copy()
$ws.Range("L21").PasteSpecial(-4122)

$ws.Range("D21").Value = "GRB"
$ws.Range("E21").Value = 8.22
$ws.Range("F21").Value = 37.9
$ws.Range("G21").Value = 0.17
$ws.Range("H21").Value = 0.2
$ws.Range("I21").Value = 2
$ws.Range("L21").Value = "2ml of amonium chloride added "

# ---- Rows 22-25: new rows without a pH (E) reading, same layout as row 19
#      but skipping column E ----
$ws.Range("A19:D19").Copy()
$ws.Range("A22:D22").PasteSpecial(-4122)
$ws.Range("F19:I19").Copy()
$ws.Range("F22:I22").PasteSpecial(-4122)
$ws.Range("A22").Value = 45521
$ws.Range("B22").Value = 21
$ws.Range("C22").Value = $instructions
$ws.Range("D22").Value = "MWW/GRB"
$ws.Range("F22").Value = 38.2
$ws.Range("G22").Value = 0.11
$ws.Range("H22").Value = 0.2
$ws.Range("I22").Value = 2

$ws.Range("A19:D19").Copy()
$ws.Range("A23:D23").PasteSpecial(-4122)
$ws.Range("F19:I19").Copy()
$ws.Range("F23:I23").PasteSpecial(-4122)
$ws.Range("A23").Value = 45522
$ws.Range("B23").Value = 22
$ws.Range("C23").Value = $instructions
$ws.Range("D23").Value = "MWW"
$ws.Range("F23").Value = 38.9
$ws.Range("G23").Value = 0.06
$ws.Range("H23").Value = 0.2
$ws.Range("I23").Value = 2

$ws.Range("A19:D19").Copy()
$ws.Range("A24:D24").PasteSpecial(-4122)
$ws.Range("F19:I19").Copy()
$ws.Range("F24:I24").PasteSpecial(-4122)
$ws.Range("A24").Value = 45523
$ws.Range("B24").Value = 23
$ws.Range("C24").Value = $instructions
$ws.Range("D24").Value = "GRB"
$ws.Range("F24").Value = 39.3
$ws.Range("G24").Value = 0.06
$ws.Range("H24").Value = 0.35
$ws.Range("I24").Value = 1

$ws.Range("A19:D19").Copy()
$ws.Range("A25:D25").PasteSpecial(-4122)
$ws.Range("F19:I19").Copy()
$ws.Range("F25:I25").PasteSpecial(-4122)
$ws.Range("A25").Value = 45524
$ws.Range("B25").Value = 24
$ws.Range("C25").Value = $instructions
$ws.Range("D25").Value = "MWW"
$ws.Range("F25").Value = 40.2
$ws.Range("G25").Value = 0.11
$ws.Range("H25").Value = 0.2
$ws.Range("I25").Value = 1

# ---- Row 26: new row with a pH (E) reading, full row 19 layout ----
$ws.Range("A19:I19").Copy()
$ws.Range("A26:I26").PasteSpecial(-4122)
$ws.Range("A26").Value = 45525
$ws.Range("B26").Value = 25
$ws.Range("C26").Value = $instructions
$ws.Range("D26").Value = "MWW/GRB"
$ws.Range("E26").Value = 8.22
$ws.Range("F26").Value = 41.8
$ws.Range("G26").Value = 0.3
$ws.Range("H26").Value = 0.2
$ws.Range("I26").Value = 2

# ---- Rows 27-31: bare rows (date / cycle day / instructions only),
#      matching the pre-edit pattern used by rows 20/21 ----
$ws.Range("A20:C20").Copy()

$ws.Range("A27:C27").PasteSpecial(-4122)
$ws.Range("A27").Value = 45526
$ws.Range("B27").Value = 26
$ws.Range("C27").Value = $instructions

$ws.Range("A28:C28").PasteSpecial(-4122)
$ws.Range("A28").Value = 45527
$ws.Range("B28").Value = 27
$ws.Range("C28").Value = $instructions

$ws.Range("A29:C29").PasteSpecial(-4122)
$ws.Range("A29").Value = 45528
$ws.Range("B29").Value = 28
$ws.Range("C29").Value = $instructions

$ws.Range("A30:C30").PasteSpecial(-4122)
$ws.Range("A30").Value = 45529
$ws.Range("B30").Value = 29
$ws.Range("C30").Value = $instructions

$ws.Range("A31:C31").PasteSpecial(-4122)
$ws.Range("A31").Value = 45530
$ws.Range("B31").Value = 30
$ws.Range("C31").Value = $instructions
